# Update the division fact-practice table: each cell's text is a unique
# "A÷B=" expression, so a straight Find/Replace per pair is unambiguous.
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
$d = $word.ActiveDocument

$d.Content.Find.Execute("21÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷8=", 2) | Out-Null
$d.Content.Find.Execute("63÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷4=", 2) | Out-Null
$d.Content.Find.Execute("51÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷5=", 2) | Out-Null
$d.Content.Find.Execute("17÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷6=", 2) | Out-Null
$d.Content.Find.Execute("70÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷2=", 2) | Out-Null
$d.Content.Find.Execute("44÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷2=", 2) | Out-Null
$d.Content.Find.Execute("84÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "38÷9=", 2) | Out-Null
$d.Content.Find.Execute("23÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷6=", 2) | Out-Null
$d.Content.Find.Execute("28÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷3=", 2) | Out-Null
$d.Content.Find.Execute("72÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=", 2) | Out-Null
$d.Content.Find.Execute("54÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷9=", 2) | Out-Null
$d.Content.Find.Execute("39÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷5=", 2) | Out-Null
$d.Content.Find.Execute("74÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=", 2) | Out-Null
$d.Content.Find.Execute("53÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷7=", 2) | Out-Null
$d.Content.Find.Execute("60÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=", 2) | Out-Null
$d.Content.Find.Execute("19÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=", 2) | Out-Null
$d.Content.Find.Execute("83÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷7=", 2) | Out-Null
$d.Content.Find.Execute("16÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷5=", 2) | Out-Null
$d.Content.Find.Execute("47÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷6=", 2) | Out-Null
$d.Content.Find.Execute("90÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷5=", 2) | Out-Null
$d.Content.Find.Execute("20÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷9=", 2) | Out-Null
$d.Content.Find.Execute("19÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷5=", 2) | Out-Null
$d.Content.Find.Execute("61÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=", 2) | Out-Null
$d.Content.Find.Execute("79÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷7=", 2) | Out-Null
$d.Content.Find.Execute("83÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷4=", 2) | Out-Null
